# Generate Report for Handoff
# - Files that were "low" priority and pending handoff are now generated
#   ("ht" priority/status marker) with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh       = $wb.Worksheets.Item("zh-cn")
$wsDe       = $wb.Worksheets.Item("de-de")

foreach ($r in 4..7) {
    # Overview: "Latest HO Xliff Generate Date" refreshed for the four
    # files that were ready for handoff.
    $wsOverview.Range("G" + $r).Value = "2016-09-03 18:43:49"

    # zh-cn: Priority flips from "low" to "ht", and the handoff datetime
    # for those rows is refreshed.
    $wsZh.Range("E" + $r).Value = "ht"
    $wsZh.Range("H" + $r).Value = "2016-09-03 18:43:44"

    # de-de: Priority flips from "low" to "ht"; its handoff datetime
    # shares the same underlying value as the Overview sheet, so it
    # also picks up the refreshed timestamp.
    $wsDe.Range("E" + $r).Value = "ht"
    $wsDe.Range("H" + $r).Value = "2016-09-03 18:43:49"
}
